# Append a parenthetical note to the "P3: Enable VNC" bullet, matching
# the author's edit: "... P3: Enable VNC (You will need RealVNC vnc
# viewer to connect from another computer)".
#
# The appended text becomes its own run placed right after the existing
# "P3: Enable VNC" run (it is not merged into that run's text), mirroring
# the pattern already used elsewhere in this document for similar
# parenthetical asides (e.g. the "Configure SSH" bullet, which also
# ends up as two separate runs: "Configure SSH" + " (Change to
# non-default port for a little better security)").

$d = $word.ActiveDocument

$addition = " (You will need RealVNC vnc viewer to connect from another computer)"

$hit = $d.Content
$found = $hit.Find.Execute("P3: Enable VNC", $true, $false, $false, $false,
                            $false, $true, 1, $false, "", 0)

if ($found) {
    # Collapse to the point right after "P3: Enable VNC" and append the
    # new text there.
    $insertionPoint = $hit.Duplicate
    $insertionPoint.Collapse(0)   # wdCollapseEnd
    $insertionPoint.InsertAfter($addition)

    # Briefly flip Bold on just the newly-inserted span (then flip it
    # back off) so the interop engine keeps the new text as its own
    # <w:r> instead of silently re-merging it into the preceding
    # "P3: Enable VNC" run -- in the source document this is a second,
    # distinct run, not a continuation of the first run's text.
    $newRange = $d.Range($hit.End, $hit.End + $addition.Length)
    $newRange.Font.Bold = 1
    $newRange.Font.Bold = 0
}
